# Auto-generated COM-interop edit script
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # Input
$ws2 = $wb.Worksheets.Item(2)   # Output

# ===== Sheet1 (Input) =====
# Delete row 25 (its content was folded into row 24; row itself is removed)
$ws1.Rows.Item(25).Delete()

# Update changed cells (rows 2-24)
$ws1.Range("D2").Value = 0.582
$ws1.Range("E2").Value = 0.334
$ws1.Range("F2").Value = 0.264
$ws1.Range("H2").Value = 0
$ws1.Range("B3").Value = "TRA_FT_COA"
$ws1.Range("C3").Value = "PRI_COA_BCO"
$ws1.Range("E3").Value = 0.105
$ws1.Range("F3").Value = 0.105
$ws1.Range("G3").Value = 0
$ws1.Range("H3").Value = 0.0730300317502966
$ws1.Range("E4").Value = 2.838766910623198
$ws1.Range("G4").Value = 1.419383455311599
$ws1.Range("D5").Value = 198982.8156135253
$ws1.Range("E5").Value = 204242.1749637213
$ws1.Range("F5").Value = 176809.0812363035
$ws1.Range("G5").Value = 229372.870751458
$ws1.Range("H5").Value = 191828.4347297169
$ws1.Range("I5").Value = 118852.1017272293
$ws1.Range("J5").Value = 20160.00890140725
$ws1.Range("K5").Value = 12337.94626304703
$ws1.Range("L5").Value = 1565.086332553496
$ws1.Range("D6").Value = 7418.568219084142
$ws1.Range("E6").Value = 14709.30790823598
$ws1.Range("F6").Value = 19335.4262232156
$ws1.Range("G6").Value = 21526.79399099144
$ws1.Range("H6").Value = 26165.95440800993
$ws1.Range("I6").Value = 28754.99863628495
$ws1.Range("J6").Value = 2926.452905042988
$ws1.Range("K6").Value = 2884.195230322684
$ws1.Range("L6").Value = 939.0517995320976
$ws1.Range("J7").Value = 60102.8394409051
$ws1.Range("K7").Value = 67948.70543292888
$ws1.Range("L7").Value = 20003.50253104742
$ws1.Range("J8").Value = 14359.12892074427
$ws1.Range("K8").Value = 12968.99408445752
$ws1.Range("L8").Value = 8794.085987936915
$ws1.Range("D9").Value = 241.4490648374262
$ws1.Range("G9").Value = 262.7649967022483
$ws1.Range("H9").Value = 611.5562900097099
$ws1.Range("I9").Value = 1007.650494806335
$ws1.Range("J9").Value = 1088.739943476629
$ws1.Range("K9").Value = 1377.531199642253
$ws1.Range("L9").Value = 1696.504881919234
$ws1.Range("D10").Value = 2.058542131979698
$ws1.Range("E10").Value = 8.076904107402811
$ws1.Range("G10").Value = 39.0685286250268
$ws1.Range("H10").Value = 168.6577354556589
$ws1.Range("I10").Value = 425.6095355729003
$ws1.Range("J10").Value = 625.8261300505758
$ws1.Range("K10").Value = 1090.700880055701
$ws1.Range("L10").Value = 1696.504881919234
$ws1.Range("D11").Value = 75161.48772542756
$ws1.Range("E11").Value = 63058.06832603392
$ws1.Range("G11").Value = 43918.95812426758
$ws1.Range("H11").Value = 27559.44845161116
$ws1.Range("I11").Value = 4877.669203698396
$ws1.Range("J11").Value = 2981.470868606867
$ws1.Range("K11").Value = 76.21218936394932
$ws1.Range("L11").Value = 94.03398386871736
$ws1.Range("D12").Value = 14059.61842364205
$ws1.Range("E12").Value = 12056.38539063793
$ws1.Range("G12").Value = 8397.083817470029
$ws1.Range("H12").Value = 5269.227880056473
$ws1.Range("I12").Value = 932.5858100152909
$ws1.Range("J12").Value = 570.0422289663444
$ws1.Range("K12").Value = 14.57138714883057
$ws1.Range("L12").Value = 17.9788245887359
$ws1.Range("D13").Value = 3276.383480154389
$ws1.Range("E13").Value = 4203.871221735594
$ws1.Range("H13").Value = 1837.296563440744
$ws1.Range("I13").Value = 325.1779469132264
$ws1.Range("J13").Value = 198.7647245737911
$ws1.Range("K13").Value = 5.080812624263288
$ws1.Range("L13").Value = 6.268932257914494
$ws1.Range("D14").Value = 179.5189976813809
$ws1.Range("E14").Value = 100.0509906883739
$ws1.Range("F14").Value = 37.3515383330611
$ws1.Range("G14").Value = 44.68282135770457
$ws1.Range("D15").Value = 9.862565402026727
$ws1.Range("E15").Value = 8.629744726773383
$ws1.Range("G15").Value = 544.6918750569678
$ws1.Range("H15").Value = 1064.423192185728
$ws1.Range("I15").Value = 3303.532621003183
$ws1.Range("J15").Value = 3118.691579497728
$ws1.Range("K15").Value = 146.2491829300208
$ws1.Range("C16").Value = "RNW_POT_BIO_GAS"
$ws1.Range("D16").Value = 0.519082389580354
$ws1.Range("E16").Value = 1.23282067525334
$ws1.Range("F16").Value = 31.83668327550674
$ws1.Range("G16").Value = 280.0472365331455
$ws1.Range("H16").Value = 980.3897822763289
$ws1.Range("I16").Value = 1225.303906516373
$ws1.Range("J16").Value = 1915.011843415477
$ws1.Range("K16").Value = 1552.032145379815
$ws1.Range("C17").Value = "HH2_BL"
$ws1.Range("D17").Value = 0
$ws1.Range("E17").Value = 0
$ws1.Range("F17").Value = 1.591834163775336
$ws1.Range("G17").Value = 15.402598009323
$ws1.Range("H17").Value = 56.02227327293309
$ws1.Range("I17").Value = 164.2583196509684
$ws1.Range("J17").Value = 228.0074721180197
$ws1.Range("K17").Value = 92.52499328225839
$ws1.Range("B18").Value = "TRA_FT_LNG"
$ws1.Range("C18").Value = "PRI_GAS_LNG"
$ws1.Range("F18").Value = 0
$ws1.Range("G18").Value = 0.9899090011670063
$ws1.Range("H18").Value = 6.152637391702759
$ws1.Range("I18").Value = 38.31239935269515
$ws1.Range("J18").Value = 225.4351405628443
$ws1.Range("K18").Value = 0
$ws1.Range("B19").Value = "TRA_FT_ETH"
$ws1.Range("C19").Value = "RNW_BIO_ETH"
$ws1.Range("D19").Value = 4.277305825242719
$ws1.Range("E19").Value = 2.138652912621359
$ws1.Range("G19").Value = 0
$ws1.Range("H19").Value = 0
$ws1.Range("I19").Value = 0
$ws1.Range("J19").Value = 0
$ws1.Range("K19").Value = 0
$ws1.Range("L19").Value = 0
$ws1.Range("B20").Value = "TRA_FT_AMM_ELCSYS_CU"
$ws1.Range("C20").Value = "ELC_CEN"
$ws1.Range("D20").Value = 0
$ws1.Range("E20").Value = 0
$ws1.Range("J20").Value = 23.00908787904872
$ws1.Range("K20").Value = 0.3947082055076144
$ws1.Range("L20").Value = 23.00908787904872
$ws1.Range("C21").Value = "HH2_WE_CU"
$ws1.Range("H21").Value = 0
$ws1.Range("I21").Value = 0
$ws1.Range("J21").Value = 186.1644382941215
$ws1.Range("K21").Value = 3.193548208197972
$ws1.Range("L21").Value = 186.1644382941215
$ws1.Range("B22").Value = "TRA_FT_AMM_ELCSYS_DT"
$ws1.Range("C22").Value = "ELC_CEN"
$ws1.Range("F22").Value = 0.172882497426718
$ws1.Range("H22").Value = 0.0211259765142098
$ws1.Range("I22").Value = 4.702092510052073
$ws1.Range("J22").Value = 4.702092510052073
$ws1.Range("K22").Value = 4.683075435335134
$ws1.Range("L22").Value = 105.6532955662002
$ws1.Range("C23").Value = "HH2_WE_DT"
$ws1.Range("F23").Value = 0
$ws1.Range("H23").Value = 0.170928355433152
$ws1.Range("I23").Value = 38.04420303587586
$ws1.Range("J23").Value = 38.04420303587586
$ws1.Range("K23").Value = 37.89033761316609
$ws1.Range("L23").Value = 854.8312095810742
$ws1.Range("B24").Value = "TRA_FT_MTH"
$ws1.Range("C24").Value = "SYN_MTH"
$ws1.Range("H24").Value = 0
$ws1.Range("I24").Value = 0.9899090011670064
$ws1.Range("J24").Value = 6.152637391702759
$ws1.Range("K24").Value = 38.31239935269517
$ws1.Range("L24").Value = 225.4351405628443

# ===== Sheet2 (Output) =====
# Insert a new row at 16 (pushes old row16 -> row17); row15 stays, new row16 is blank then filled
$ws2.Rows.Item(16).Insert()

# Update changed cells (rows 1-15, unaffected by the insert since insert was at row16)
$ws2.Range("D2").Value = 0.582
$ws2.Range("H2").Value = 0.0730300317502966
$ws2.Range("D4").Value = 8830.748095487703
$ws2.Range("E4").Value = 9326.65791604474
$ws2.Range("F4").Value = 9066.4668709098
$ws2.Range("G4").Value = 11467.44318411311
$ws2.Range("H4").Value = 10073.28164015216
$ws2.Range("I4").Value = 7031.09395749137
$ws2.Range("J4").Value = 4251.485748159676
$ws2.Range("K4").Value = 4197.625691597962
$ws2.Range("L4").Value = 1369.137523717798
$ws2.Range("G5").Value = 280.1727241433552
$ws2.Range("H5").Value = 725.9144169418577
$ws2.Range("I5").Value = 1336.405766490116
$ws2.Range("J5").Value = 1601.619271263929
$ws2.Range("K5").Value = 2310.382195722
$ws2.Range("L5").Value = 3180.946653598565
$ws2.Range("D6").Value = 3826.552589408666
$ws2.Range("G6").Value = 2269.864341396109
$ws2.Range("H6").Value = 1424.355494314213
$ws2.Range("I6").Value = 252.0926694862386
$ws2.Range("J6").Value = 154.091415156376
$ws2.Range("K6").Value = 3.938876020862829
$ws2.Range("L6").Value = 4.859960162513031
$ws2.Range("D7").Value = 1605.070043182009
$ws2.Range("E7").Value = 1981.395153867322
$ws2.Range("G7").Value = 383.3687658878648
$ws2.Range("L7").Value = 413.9949531830487
$ws2.Range("D8").Value = 179.5189976813809
$ws2.Range("E8").Value = 100.0509906883739
$ws2.Range("F8").Value = 37.3515383330611
$ws2.Range("G8").Value = 44.68282135770457
$ws2.Range("D9").Value = 10.29859460927422
$ws2.Range("E9").Value = 9.783664878810509
$ws2.Range("G9").Value = 834.7605630042511
$ws2.Range("H9").Value = 2088.074685242437
$ws2.Range("I9").Value = 4666.260300542878
$ws2.Range("J9").Value = 5232.035090353493
$ws2.Range("K9").Value = 1781.109308481964
$ws2.Range("G10").Value = 0.9899090011670063
$ws2.Range("K10").Value = 0
$ws2.Range("L10").Value = 0
$ws2.Range("H12").Value = 0
$ws2.Range("I12").Value = 0
$ws2.Range("J12").Value = 96.42899556583146
$ws2.Range("K12").Value = 1.654186206718276
$ws2.Range("L12").Value = 96.42899556583146
$ws2.Range("H13").Value = 0.08853704702773381
$ws2.Range("I13").Value = 19.70604224667278
$ws2.Range("J13").Value = 19.70604224667278
$ws2.Range("K13").Value = 19.62634341535906
$ws2.Range("L13").Value = 442.7833568728934
$ws2.Range("K14").Value = 38.31239935269517
$ws2.Range("H15").Value = 0
$ws2.Range("I15").Value = 0.9379261261528213
$ws2.Range("J15").Value = 5.88872210961464
$ws2.Range("K15").Value = 36.71525896327071
$ws2.Range("L15").Value = 216.1033327415355

# Fill brand-new row 16 (HH2_DEL_TRA_LH2_C_2_NEW) entirely
$ws2.Range("A16").Value = 14
$ws2.Range("B16").Value = "HH2_DEL_TRA_LH2_C_2_NEW"
$ws2.Range("C16").Value = "TRA_LH2"
$ws2.Range("D16").Value = 0
$ws2.Range("E16").Value = 0
$ws2.Range("F16").Value = 0
$ws2.Range("G16").Value = 0
$ws2.Range("H16").Value = 0.00870511425462459
$ws2.Range("I16").Value = 0.008705114254624557
$ws2.Range("J16").Value = 0.008705114254624471
$ws2.Range("K16").Value = 0.008705114254624592
$ws2.Range("L16").Value = 0

# Row 17 = old row16 shifted down; update the cells that changed
$ws2.Range("A17").Value = 15
$ws2.Range("H17").Value = 0
$ws2.Range("I17").Value = 22.54219698339416
$ws2.Range("J17").Value = 633.6580776784569
$ws2.Range("K17").Value = 1191.55474771288
$ws2.Range("L17").Value = 1339.118844065566
